$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

$ws.Unprotect()

$ws.Range("A6").Value  = "- Antibody label: the COVIC label for the antibody (e.g. COVIC 1)"
$ws.Range("A7").Value  = "- Tested antigen: the name of the tested antigen (e.g. Spike protein 1)"
$ws.Range("A8").Value  = "- n: the number of runs for the assay (e.g. 6)"
$ws.Range("A9").Value  = "- on rate; Ka in M^-1s^-1: an SPR assay measuring on rate [Ka] in M^-1s^-1 (e.g. 491000)"
$ws.Range("A10").Value = "- Standard deviation in M^-1s^-1: The standard deviation of the value in 'Standard deviation in M^-1s^-1'"
$ws.Range("A11").Value = "- off rate; Kd in 1/s: an SPR assay measuring off rate [Kd] in 1/s (e.g. 126)"
$ws.Range("A12").Value = "- Standard deviation in 1/s: The standard deviation of the value in 'Standard deviation in 1/s'"
$ws.Range("A13").Value = "- dissociation constant; KD in nM: an SPR assay measuring dissociation constant [KD] in nM (e.g. <.1)"
$ws.Range("A14").Value = "- Standard deviation in nM: The standard deviation of the value in 'Standard deviation in nM'"
$ws.Range("A15").Value = "- Qualitiative measure: the qualitative measure of the assay (e.g. positive)"
$ws.Range("A16").Value = "- Comment: general comments about the assay (e.g. did not bind positive control)"

$ws.Protect($null, $false, $false, $false, $true, $true, $true, $true, $true, $true, $true, $true, $true, $false, $false, $true, $false)

